$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posibles")
Write-Output $ws.Name
